$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Top "Bad Drivers" section ---
$ws.Range("C3").Value = 589
$ws.Range("D3").Value = 86.5

$ws.Range("C5").Value = 367
$ws.Range("D5").Value = 98.5

$ws.Range("C6").Value = 1261

# --- Bottom "Good Drivers" section (rows 14-29) ---
# row, A (name), B (count), D (pct), E (vintage date, or $null to clear)
$rows = @(
    @(14, "Intel(R) Wi-Fi 6E AX210 160MHz - 22.0.1.5", 156943, 100, $null),
    @(15, "Intel(R) Wi-Fi 6E AX210 160MHz - 23.120.0.3", 34181, 99.90000000000001, "2025-02-05"),
    @(16, "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4", 445055, 99.90000000000001, "2024-11-10"),
    @(17, "MediaTek Wi-Fi 6E MT7922 (RZ616) 160MHz PCIe Adapter - 3.4.0.1088", 86276, 99.90000000000001, "2024-08-07"),
    @(18, "MediaTek Wi-Fi 6E MT7922 (RZ616) 160MHz PCIe Adapter - 3.3.0.1030", 17891, 100, "2024-05-09"),
    @(19, "Intel(R) Wi-Fi 6E AX210 160MHz - 23.20.1.1", 13533, 100, "2023-12-19"),
    @(20, "Intel(R) Wi-Fi 6E AX210 160MHz - 22.170.2.1", 19083, 100, "2022-08-30"),
    @(21, "Intel(R) Wi-Fi 6E AX210 160MHz - 22.100.0.3", 12988, 100, "2022-05-01"),
    @(22, "Intel(R) Wi-Fi 6E AX210 160MHz - 22.130.0.5", 18738, 99.90000000000001, "2022-03-14"),
    @(23, "Intel(R) Wi-Fi 6E AX210 160MHz - 22.110.1.1", 42024, 100, "2022-01-01"),
    @(24, "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9", 77849, 99.90000000000001, "2021-08-18"),
    @(25, "Intel(R) Wi-Fi 6E AX210 160MHz - 22.70.0.6", 15504, 100, "2021-06-28"),
    @(26, "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1", 34244, 100, "2021-04-27"),
    @(27, "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2", 59673, 100, "2020-08-05"),
    @(28, "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6", 113652, 100, "2020-01-06"),
    @(29, "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1", 56018, 100, "2019-12-14")
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $name = $r[1]
    $count = $r[2]
    $pct = $r[3]
    $vintage = $r[4]

    $ws.Range("A$rowNum").Value = $name
    $ws.Range("B$rowNum").Value = $count
    $ws.Range("D$rowNum").Value = $pct

    if ($null -eq $vintage) {
        $ws.Range("E$rowNum").Value = ""
    } else {
        # Force text (not auto-parsed date) so the stored cell matches the
        # literal "yyyy-mm-dd" inline string the driver-vintage column uses.
        $ws.Range("E$rowNum").NumberFormat = "@"
        $ws.Range("E$rowNum").Value = $vintage
    }
}
